$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 597, shifting existing rows 597-699 down to 599-701
$ws.Rows.Item(597).Resize(2).Insert()

# Populate new row 597 with data
$ws.Cells.Item(597, 1).Value = 7
$ws.Cells.Item(597, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(597, 3).Value = "Ñuble"
$ws.Cells.Item(597, 4).Value = 44694
$ws.Cells.Item(597, 5).Value = 16
$ws.Cells.Item(597, 6).Value = "Fruta"
$ws.Cells.Item(597, 7).Value = 100102
$ws.Cells.Item(597, 8).Value = "Cítricos"
$ws.Cells.Item(597, 9).Value = 100102003
$ws.Cells.Item(597, 10).Value = "Limón"
$ws.Cells.Item(597, 11).Value = "Sin especificar"
$ws.Cells.Item(597, 12).Value = "1a amarillo"
$ws.Cells.Item(597, 13).Value = 120
$ws.Cells.Item(597, 14).Value = 14000
$ws.Cells.Item(597, 15).Value = 15000
$ws.Cells.Item(597, 16).Value = 14500
$ws.Cells.Item(597, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(597, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(597, 19).Value = 906
$ws.Cells.Item(597, 20).Value = 16

# Populate new row 598 with data
$ws.Cells.Item(598, 1).Value = 7
$ws.Cells.Item(598, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(598, 3).Value = "Ñuble"
$ws.Cells.Item(598, 4).Value = 44694
$ws.Cells.Item(598, 5).Value = 16
$ws.Cells.Item(598, 6).Value = "Fruta"
$ws.Cells.Item(598, 7).Value = 100102
$ws.Cells.Item(598, 8).Value = "Cítricos"
$ws.Cells.Item(598, 9).Value = 100102003
$ws.Cells.Item(598, 10).Value = "Limón"
$ws.Cells.Item(598, 11).Value = "Sin especificar"
$ws.Cells.Item(598, 12).Value = "2a amarillo"
$ws.Cells.Item(598, 13).Value = 120
$ws.Cells.Item(598, 14).Value = 12000
$ws.Cells.Item(598, 15).Value = 13000
$ws.Cells.Item(598, 16).Value = 12500
$ws.Cells.Item(598, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(598, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(598, 19).Value = 781
$ws.Cells.Item(598, 20).Value = 16
